# Add "Save" column (H) to the s_vals worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - copy the header formatting from the neighboring "sum" column (G1)
# and then set the new label.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Determine save flag per row based on the "sum" column (G): 1 if sum > 10, else 0.
for ($r = 2; $r -le 68; $r++) {
    $g = $ws.Cells.Item($r, 7).Value()
    if ($g -gt 10) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
